# Regenerate merged AHB files
# - Rename header cells in row 1 from *_old / *_new suffix to *_FV2304 / *_FV2310
# - Add an Excel Table (ListObject) over A1:U58
# - Freeze the header row (pane split) and restore a selection in the frozen pane

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2304Suffix = "_FV2304"
$fv2310Suffix = "_FV2310"

# Columns A-J use the "_old" -> "_FV2304" suffix, column K is "diff" (unchanged),
# columns L-U use the "_new" -> "_FV2310" suffix.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = [string]$cell.Value2
    if ($text -like "*$oldSuffix") {
        $base = $text.Substring(0, $text.Length - $oldSuffix.Length)
        $cell.Value2 = $base + $fv2304Suffix
    }
    elseif ($text -like "*$newSuffix") {
        $base = $text.Substring(0, $text.Length - $newSuffix.Length)
        $cell.Value2 = $base + $fv2310Suffix
    }
}

# Create the Excel table (ListObject) covering the full used range A1:U58
$tableRange = $ws.Range("A1:U58")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

# Freeze panes at the header row (row 2 is the top-left cell of the scrollable area)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
